# Fresh rolls, minor display updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the parenthetical note text (C1 / shared string)
$ws.Range("C1").Value = "(Rolling on the flat wood desk, pick up and reroll after some shaking, avoid bumping into each other or objects near end of roll)"

# New dice rolls: Orange (col A) / Blue (col B), appended as rows 2303-2402
$newRolls = @(
    @(6, 14),
    @(2, 19),
    @(2, 4),
    @(3, 8),
    @(3, 4),
    @(18, 4),
    @(3, 8),
    @(12, 12),
    @(9, 8),
    @(2, 13),
    @(11, 14),
    @(9, 19),
    @(19, 4),
    @(4, 5),
    @(17, 13),
    @(6, 4),
    @(12, 8),
    @(8, 10),
    @(3, 5),
    @(1, 10),
    @(6, 20),
    @(3, 7),
    @(3, 14),
    @(17, 20),
    @(4, 15),
    @(10, 4),
    @(2, 14),
    @(10, 2),
    @(10, 4),
    @(18, 20),
    @(7, 19),
    @(9, 17),
    @(1, 1),
    @(4, 20),
    @(18, 2),
    @(3, 19),
    @(13, 20),
    @(12, 17),
    @(9, 14),
    @(13, 8),
    @(14, 6),
    @(17, 7),
    @(13, 5),
    @(11, 7),
    @(1, 8),
    @(19, 14),
    @(6, 14),
    @(3, 12),
    @(2, 6),
    @(18, 18),
    @(17, 9),
    @(17, 18),
    @(20, 14),
    @(9, 16),
    @(3, 1),
    @(13, 7),
    @(18, 18),
    @(10, 1),
    @(13, 15),
    @(7, 7),
    @(16, 15),
    @(6, 14),
    @(12, 20),
    @(12, 19),
    @(11, 5),
    @(10, 12),
    @(2, 18),
    @(4, 9),
    @(6, 8),
    @(20, 20),
    @(18, 5),
    @(13, 8),
    @(3, 3),
    @(12, 13),
    @(12, 5),
    @(15, 20),
    @(13, 7),
    @(7, 4),
    @(9, 9),
    @(1, 10),
    @(4, 8),
    @(9, 1),
    @(15, 17),
    @(19, 6),
    @(2, 19),
    @(10, 9),
    @(3, 19),
    @(11, 9),
    @(4, 13),
    @(3, 7),
    @(11, 10),
    @(15, 19),
    @(10, 10),
    @(11, 7),
    @(16, 12),
    @(15, 8),
    @(18, 14),
    @(14, 15),
    @(14, 12),
    @(2, 9)
)

$startRow = 2303
for ($i = 0; $i -lt $newRolls.Count; $i++) {
    $r = $startRow + $i
    $pair = $newRolls[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

$lastRow = $startRow + $newRolls.Count

# Keep the frozen-pane view scrolled near the bottom of the data, mirroring
# Excel's own behavior of following new entries down the sheet.
try {
    $excel.ActiveWindow.ScrollRow = $lastRow - 33
} catch {}

$ws.Range("A" + $lastRow).Select()
